$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add latitude/longitude coordinates for Roatán habitat sites (Honduras)
$ws.Range("C12").Value = 16.393
$ws.Range("D12").Value = -86.274

$ws.Range("C13").Value = 16.398
$ws.Range("D13").Value = -86.269

$ws.Range("C19").Value = 16.358
$ws.Range("D19").Value = -86.289

$ws.Range("C20").Value = 16.362
$ws.Range("D20").Value = -86.279

$ws.Range("C21").Value = 16.374
$ws.Range("D21").Value = -86.283
